# Add the TrappedWater solution link to column E (row 2), resize column E /
# row 2 to fit the new content, and move the sheet's active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$solutionUrl = "https://github.com/Gershon-Tadepalli/DS-Algo/blob/master/DS-AlgoPractice/DS-AlgoPractice/TrapWater.cs"

# 1) Put the link text into E2 and turn it into a real hyperlink (mirrors the
#    existing QuestionLink hyperlink already present on C2).
$ws.Range("E2").Value = $solutionUrl
$ws.Hyperlinks.Add($ws.Range("E2"), $solutionUrl)

# Excel's Hyperlinks.Add always (re)stamps its own "Hyperlink" cell style;
# put E2 back on the same style already used by C2 so no redundant style is
# introduced.
$ws.Range("E2").Style = $ws.Range("C2").Style

# 2) Column E needs to be wide enough for the new text (raw stored width of
#    29 characters). 28.16 is comfortably inside the COM ColumnWidth band
#    that rounds to that stored width.
$ws.Columns("E").ColumnWidth = 28.16

# 3) Row 2 grows to fit the wrapped text.
$ws.Rows(2).RowHeight = 64.5

# 4) The sheet's selection moves to E5.
[void]$ws.Range("E5").Select()
